# Moves the four SWOT "Content Placeholder" boxes (Strengths/Weaknesses on
# the top row, Opportunities/Threats on the bottom row) up by 198780 EMU on
# the "Surtitle, Title Only" slide layout, matching the target OOXML:
#   id=14 (idx 11) y: 1373188 -> 1174408
#   id=15 (idx 12) y: 1373188 -> 1174408
#   id=16 (idx 13) y: 4091188 -> 3892408
#   id=17 (idx 14) y: 4091188 -> 3892408
# x offsets, widths and heights are untouched.

$EMU_PER_POINT = 12700

# PowerPoint's Shape.Top/.Left are expressed in points (single-precision
# under the hood), so converting an exact EMU target straight to points and
# back can truncate to one EMU short. A tiny nudge (well under a hundredth
# of a point) keeps the rounding on the correct side without perceptibly
# moving the shape.
function Set-TopFromEmu($Shape, $Emu) {
    $Shape.Top = ($Emu / $EMU_PER_POINT) + 0.00001
}

$p = $ppt.ActivePresentation
$master = $p.Slides.Item(1).Master
$layout = $master.CustomLayouts.Item(1)

$targets = @{
    14 = 1174408
    15 = 1174408
    16 = 3892408
    17 = 3892408
}

for ($i = 1; $i -le $layout.Shapes.Count; $i++) {
    $shp = $layout.Shapes.Item($i)
    if ($targets.ContainsKey($shp.Id)) {
        Set-TopFromEmu $shp $targets[$shp.Id]
    }
}
